# Update the shared-string description for the CO2-HTC flow (column F header,
# row 1) on every sheet: the unit changes from "kg" to "kilogram" and gains an
# explicit compartment tuple. The leading apostrophe must be doubled so the
# COM layer doesn't treat it as an Excel "force text" prefix marker and strip
# it from the stored string.
$newDescription = "''Carbon dioxide, non-fossil, resource correction' (kilogram, None, ('natural resource', 'in air'))"

# New column-H (CO2 - HTC) values for rows 4-14, expressed as plain decimals
# (the COM layer's expression parser chokes on scientific-notation literals).
$col_h = @{
    4  = 0.00000103
    5  = 0.0000559
    6  = 0.000000533
    7  = 0.000199
    8  = 0.0000775
    9  = 0.000000113
    10 = 0.000000000000679
    11 = 0.000000000017
    12 = 0.00000000000142
    13 = 0.000000311
    14 = 0.00000597
}

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count()
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $ws.Range("F1").Value = $newDescription

    foreach ($row in $col_h.Keys) {
        $ws.Cells.Item($row, 8).Value = $col_h[$row]
    }
}
